$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109. This shifts the existing rows 109-159
# down to 110-160 (preserving all their values/formatting), matching the
# target diff where every record from the old row N moved to row N+1.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record's data.
$ws.Cells.Item(109, 1).Value  = 5
$ws.Cells.Item(109, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value  = "Maule"
$ws.Cells.Item(109, 4).Value  = 44455
$ws.Cells.Item(109, 5).Value  = 7
$ws.Cells.Item(109, 6).Value  = 100112003
$ws.Cells.Item(109, 7).Value  = "Ajo"
$ws.Cells.Item(109, 8).Value  = "Chino"
$ws.Cells.Item(109, 9).Value  = "Primera"
$ws.Cells.Item(109, 10).Value = 200
$ws.Cells.Item(109, 11).Value = 17000
$ws.Cells.Item(109, 12).Value = 17000
$ws.Cells.Item(109, 13).Value = 17000
$ws.Cells.Item(109, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(109, 15).Value = "China"
$ws.Cells.Item(109, 16).Value = 1700
$ws.Cells.Item(109, 17).Value = 10
$ws.Cells.Item(109, 18).Value = "Hortaliza"
